# Updated cryptos list with GitHub Actions
# Applies per-cell updates to match the target snapshot of cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.480.64"
$ws.Range("E2").Value = "  -1.10%  "

# Row 3
$ws.Range("D3").Value = "1.920.95"
$ws.Range("E3").Value = "  +1.49%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4693"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.23%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2867"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06743"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.46%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "106.04"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.30"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.90%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07743"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.09%  "

# Row 13
$ws.Range("D13").Value = "1.903.05"
$ws.Range("E13").Value = "  +0.85%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.296"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.39%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6575"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "289.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.21%  "

# Row 17
$ws.Range("D17").Value = "30.485.91"
$ws.Range("E17").Value = "  -0.98%  "

# Row 18
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007581"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.60%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.93"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.02%  "

# Row 21
$ws.Range("D21").Value = "2.148.13"
$ws.Range("E21").Value = "  +1.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.249"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.190"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.361"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.23%  "

# Row 26
$ws.Range("E26").Value = "  +0.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.124"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1066"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -8.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.367"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.171"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.58%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.977"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05030"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7386"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.153"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02093"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.84%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.03%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.682"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.99%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.055"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.45%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.26"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8689"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.29%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.854"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.43%  "

# Row 43
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4240"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.56%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.22"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "49.92"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +16.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.186"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.299"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.28%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.96"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.33%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1213"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.86%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2471"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +10.69%  "
